{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Diff being applied: the first paragraph, \"This is a Microsoft word\n// document.\", gets \" (Changed main)\" appended after the existing text,\n// added as three separate runs: \" (\", \"Changed main\", \")\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the target paragraph robustly by its text rather than assuming index 0.\nlet target = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"This is a Microsoft word document.\") !== -1) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  target = paragraphs.items[0];\n}\n\n// Append the new text at the end of the paragraph, as three distinct runs.\ntarget.insertText(\" (\", Word.InsertLocation.end);\nawait context.sync();\n\ntarget.insertText(\"Changed main\", Word.InsertLocation.end);\nawait context.sync();\n\ntarget.insertText(\")\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Diff being applied: the first paragraph, \"This is a Microsoft word\n# document.\", gets \" (Changed main)\" appended after the existing text,\n# added as three separate runs: \" (\", \"Changed main\", \")\".\n\n$d = $word.ActiveDocument\n\n# Locate the target text robustly via Find, then collapse to its end so the\n# insertions land right after the existing sentence (and stay inside the\n# same paragraph, not spilling into the next one).\n$rng = $d.Content\n$rng.Find.Execute(\"This is a Microsoft word document.\")\n$rng.Collapse(0)   # wdCollapseEnd\n\n$rng.InsertAfter(\" (\")\n$rng.Collapse(0)\n\n$rng.InsertAfter(\"Changed main\")\n$rng.Collapse(0)\n\n$rng.InsertAfter(\")\")\n$rng.Collapse(0)\n"}
